# Apply the commit "Deploy the implementation guide." changes:
# 1. Rename the second sheet tab.
# 2. Update the Metadata sheet: Date, Contact values, and insert a new
#    "Jurisdiction" row (pushing Description/Purpose/Copyright/Immutable down).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)
$inc  = $wb.Worksheets.Item(2)

# --- Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0" ---
$inc.Name = "Include #0"

# --- Update Date value (row 8) ---
$meta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# --- Update Contact value (row 10) ---
$meta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# --- Insert a new "Jurisdiction" row after "Contact" (new row 11) ---
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Match the formatting of the other data rows (border/alignment style)
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
